# Applies the "cryptos list" data refresh described by the commit
# "Updated cryptos list on Sun Jun  9 03:30:15 UTC 2024 with GitHub Actions"
#
# Every write goes through the same Set-TextValue helper so that values
# which *look* numeric (e.g. "672.90") are stored as literal text,
# matching the inlineStr / shared-string cells already used throughout
# this sheet, instead of being auto-coerced to the Number type by Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}


# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "69.174.12"
Set-TextValue $ws.Range("E2") "  -0.42%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "3.670.63"
Set-TextValue $ws.Range("E3") "  -0.47%  "

# Row 4 - TetherUSD
Set-TextValue $ws.Range("E4") "  +0.09%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "672.90"
Set-TextValue $ws.Range("E5") "  -1.24%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "156.95"
Set-TextValue $ws.Range("E6") "  -3.51%  "

# Row 7 - USDC
Set-TextValue $ws.Range("E7") "  +0.04%  "

# Row 8 - XRP
Set-TextValue $ws.Range("E8") "  -1.57%  "

# Row 9 - Dogecoin
Set-TextValue $ws.Range("E9") "  -2.08%  "

# Row 10 - Toncoin
Set-TextValue $ws.Range("D10") "6.94"
Set-TextValue $ws.Range("E10") "  -5.88%  "

# Row 11 - Cardano
Set-TextValue $ws.Range("E11") "  -2.69%  "

# Row 12 - ShibaInu
Set-TextValue $ws.Range("E12") "  -3.97%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D13") "4.288.09"
Set-TextValue $ws.Range("E13") "  -0.55%  "

# Row 14 - Avalanche
Set-TextValue $ws.Range("D14") "32.11"
Set-TextValue $ws.Range("E14") "  -4.48%  "

# Row 15 - WrappedEther
Set-TextValue $ws.Range("D15") "3.670.14"
Set-TextValue $ws.Range("E15") "  -0.57%  "

# Row 16 - WrappedBTC
Set-TextValue $ws.Range("D16") "69.161.03"
Set-TextValue $ws.Range("E16") "  -0.45%  "

# Row 17 - TRON
Set-TextValue $ws.Range("E17") "  +0.84%  "

# Row 18 - Chainlink
Set-TextValue $ws.Range("D18") "15.98"
Set-TextValue $ws.Range("E18") "  -1.68%  "

# Row 19 - Polkadot
Set-TextValue $ws.Range("D19") "6.40"
Set-TextValue $ws.Range("E19") "  -3.61%  "

# Row 20 - BitcoinCash
Set-TextValue $ws.Range("D20") "466.68"
Set-TextValue $ws.Range("E20") "  -3.61%  "

# Row 21 - Uniswap
Set-TextValue $ws.Range("D21") "9.91"
Set-TextValue $ws.Range("E21") "  -0.05%  "

# Row 22 - Polygon
Set-TextValue $ws.Range("E22") "  -3.33%  "

# Row 23 - Litecoin
Set-TextValue $ws.Range("D23") "79.58"
Set-TextValue $ws.Range("E23") "  -0.90%  "

# Row 24 - WrappedeETH
Set-TextValue $ws.Range("D24") "3.815.72"
Set-TextValue $ws.Range("E24") "  -0.45%  "

# Row 25 - Dai
Set-TextValue $ws.Range("E25") "  -0.06%  "

# Row 26 - PEPE <-> InternetComputer(DFINITY) swap (now InternetComputer)
Set-TextValue $ws.Range("B26") "InternetComputer(DFINITY)"
Set-TextValue $ws.Range("C26") "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D26") "10.86"
Set-TextValue $ws.Range("E26") "  -5.55%  "

# Row 27 - PEPE <-> InternetComputer(DFINITY) swap (now PEPE)
Set-TextValue $ws.Range("B27") "PEPE"
Set-TextValue $ws.Range("C27") "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue $ws.Range("D27") "0.0000120"
Set-TextValue $ws.Range("E27") "  -8.17%  "

# Row 28 - RenderToken
Set-TextValue $ws.Range("E28") "  -6.13%  "

# Row 29 - PancakeSwap
Set-TextValue $ws.Range("E29") "  -2.73%  "

# Row 30 - Fetch.AI
Set-TextValue $ws.Range("E30") "  -6.88%  "

# Row 31 - NEARProtocol
Set-TextValue $ws.Range("E31") "  -4.15%  "

# Row 32 - Binance-PegBSC-USD
Set-TextValue $ws.Range("D32") "0.999"
Set-TextValue $ws.Range("E32") "  -0.11%  "

# Row 33 - EthereumClassic
Set-TextValue $ws.Range("D33") "26.78"
Set-TextValue $ws.Range("E33") "  -1.26%  "

# Row 34 - ImmutableX
Set-TextValue $ws.Range("E34") "  -5.89%  "

# Row 35 - RenzoRestakedETH
Set-TextValue $ws.Range("D35") "3.663.29"
Set-TextValue $ws.Range("E35") "  +0.14%  "

# Row 36 - Kaspa
Set-TextValue $ws.Range("D36") "0.159"
Set-TextValue $ws.Range("E36") "  -4.68%  "

# Row 37 - Aptos
Set-TextValue $ws.Range("D37") "8.08"
Set-TextValue $ws.Range("E37") "  -5.12%  "

# Row 38 - Filecoin
Set-TextValue $ws.Range("D38") "6.14"
Set-TextValue $ws.Range("E38") "  -3.72%  "

# Row 40 - FirstDigitalUSD
Set-TextValue $ws.Range("E40") "  +0.00%  "

# Row 41 - Stacks
Set-TextValue $ws.Range("E41") "  -2.16%  "

# Row 42 - Monero
Set-TextValue $ws.Range("D42") "173.86"
Set-TextValue $ws.Range("E42") "  +8.07%  "

# Row 43 - Hedera
Set-TextValue $ws.Range("D43") "0.0894"
Set-TextValue $ws.Range("E43") "  -4.81%  "

# Row 44 - Mantle
Set-TextValue $ws.Range("E44") "  -1.87%  "

# Row 45 - OKB
Set-TextValue $ws.Range("D45") "47.55"
Set-TextValue $ws.Range("E45") "  -1.57%  "

# Row 46 - dogwifhat
Set-TextValue $ws.Range("E46") "  -6.65%  "

# Row 47 - FLOKI
Set-TextValue $ws.Range("D47") "0.000274"
Set-TextValue $ws.Range("E47") "  -5.47%  "

# Row 48 - ONDO <-> InjectiveProtocol swap (now InjectiveProtocol)
Set-TextValue $ws.Range("B48") "InjectiveProtocol"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D48") "27.49"
Set-TextValue $ws.Range("E48") "  -9.06%  "

# Row 49 - ONDO <-> InjectiveProtocol swap (now ONDO)
Set-TextValue $ws.Range("B49") "ONDO"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-TextValue $ws.Range("D49") "1.27"
Set-TextValue $ws.Range("E49") "  -6.36%  "

# Row 50 - SuiNetwork
Set-TextValue $ws.Range("E50") "  -4.20%  "

# Row 51 - Cosmos
Set-TextValue $ws.Range("D51") "7.75"
Set-TextValue $ws.Range("E51") "  -3.71%  "

